$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows whose data was repulled
$ws.Range("F2").Value = -10
$ws.Range("F3").Value = -7
$ws.Range("F4").Value = -13
$ws.Range("F5").Value = -9
$ws.Range("F7").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("F10").Value = -6
